# Refresh cryptocurrency price / 1h-volume snapshot (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '67.771.48'
$ws.Cells.Item(2, 5).Value = '  -0.23%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.805.71'
$ws.Cells.Item(3, 5).Value = '  +0.39%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.00'
$ws.Cells.Item(4, 5).Value = '  +0.36%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''597.49'
$ws.Cells.Item(5, 5).Value = '  +0.40%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''167.49'
$ws.Cells.Item(6, 5).Value = '  +0.40%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.12%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.14%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.74%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''6.30'
$ws.Cells.Item(10, 5).Value = '  -1.00%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.450'
$ws.Cells.Item(11, 5).Value = '  +0.04%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -1.21%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''35.99'
$ws.Cells.Item(13, 5).Value = '  -0.47%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '4.439.06'

# Row 15
$ws.Cells.Item(15, 4).Value = '3.813.66'
$ws.Cells.Item(15, 5).Value = '  +0.20%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''18.48'
$ws.Cells.Item(16, 5).Value = '  +1.28%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '67.842.41'
$ws.Cells.Item(17, 5).Value = '  +0.08%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +1.09%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.48%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''461.73'
$ws.Cells.Item(20, 5).Value = '  +0.25%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -3.61%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.40%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.72%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''83.31'
$ws.Cells.Item(24, 5).Value = '  -0.42%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''12.08'
$ws.Cells.Item(25, 5).Value = '  +1.49%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -1.94%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -0.07%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.99%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '3.953.34'
$ws.Cells.Item(29, 5).Value = '  +0.44%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''2.77'
$ws.Cells.Item(30, 5).Value = '  -0.68%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''7.38'
$ws.Cells.Item(31, 5).Value = '  +0.90%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +1.46%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''29.53'
$ws.Cells.Item(33, 5).Value = '  -1.36%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(34, 4).Value = '''0.996'
$ws.Cells.Item(34, 5).Value = '  -0.21%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Aptos'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(35, 4).Value = '''9.04'
$ws.Cells.Item(35, 5).Value = '  -1.46%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'RenzoRestakedETH'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Cells.Item(36, 4).Value = '3.745.72'
$ws.Cells.Item(36, 5).Value = '  +0.12%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.27%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''3.41'
$ws.Cells.Item(38, 5).Value = '  +1.75%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.24%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''5.78'
$ws.Cells.Item(41, 5).Value = '  +0.38%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''1.00'
$ws.Cells.Item(42, 5).Value = '  +0.16%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''48.08'
$ws.Cells.Item(44, 5).Value = '  +2.13%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''0.300'
$ws.Cells.Item(45, 5).Value = '  +0.66%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''42.72'
$ws.Cells.Item(46, 5).Value = '  -3.36%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(47, 4).Value = '''27.65'
$ws.Cells.Item(47, 5).Value = '  +9.88%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''8.34'
$ws.Cells.Item(48, 5).Value = '  -0.62%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Monero'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(49, 4).Value = '''147.87'
$ws.Cells.Item(49, 5).Value = '  -0.04%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'ONDO'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(50, 4).Value = '''1.36'
$ws.Cells.Item(50, 5).Value = '  +8.18%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Stacks'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(51, 4).Value = '''1.84'
$ws.Cells.Item(51, 5).Value = '  +0.31%  '
